$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text so numeric-looking
# strings (e.g. "1.000", "29.315.95") are stored verbatim as inline strings
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.336.70"
$ws.Range("D3").Value = "1.869.58"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "0.7264"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").Value = "241.21"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.07886"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "0.3099"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "25.29"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "0.08264"
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("D12").Value = "1.881.92"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "0.7241"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "5.251"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "90.85"
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("D16").Value = "29.393.29"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "5.869"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "244.44"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").Value = "0.000007835"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "13.25"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "2.113.39"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "8.009"
$ws.Range("E23").Value = "  +6.35%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "0.1599"
$ws.Range("E25").Value = "  +12.64%  "
$ws.Range("D26").Value = "162.63"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "8.970"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "18.29"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "4.403"
$ws.Range("E31").Value = "  +2.32%  "
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "0.05230"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "1.936"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("D35").Value = "1.188"
$ws.Range("E35").Value = "  +0.75%  "
$ws.Range("D36").Value = "0.7293"
$ws.Range("E36").Value = "  +3.03%  "
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "0.01862"
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").Value = "2.700"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "1.169.64"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "0.9047"
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "6.118"
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "72.78"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "102.07"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.017.02"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "0.5281"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "1.785"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").Value = "2.904"
$ws.Range("E50").Value = "  +5.75%  "
$ws.Range("D51").Value = "9.271"
$ws.Range("E51").Value = "  +1.24%  "

# Restore the default cell style so no stray number-format styling is
# left behind on cells that did not have one originally.
$dataRange.Style = "Normal"

